$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL): plain string assignment ---
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

# --- Numeric-looking text columns (Price / Volume% / Hora): force text type ---
# Prefix with an apostrophe so Excel stores these as text (matching the
# original inlineStr cells) instead of converting to a number/percentage,
# then reset the style back to Normal so no stray NumberFormat style gets
# attached to the cell (keeps styles.xml / cell s= identical to source).
$ws.Range("D2").Value = "'313.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.01%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'7"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'40.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.36%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'7"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.141"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.93%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'7"
$ws.Range("G4").Style = "Normal"
$ws.Range("E5").Value = "'-0.95%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'7"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'4.332"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.51%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'7"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'1.680"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.79%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'7"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.9303"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.62%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'7"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.1199"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.85%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'7"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.1816"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.48%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'7"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.09005"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.74%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'7"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.04148"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.44%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'7"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.29%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'7"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.001283"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.67%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'7"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.005834"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.17%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'7"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.007522"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.18%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'7"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'3.335"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.33%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'7"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'2.424"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.78%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'7"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.3345"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.28%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'7"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'7.614"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.22%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'7"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.1352"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.30%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'7"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.2839"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.57%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'7"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.03974"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.68%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'7"
$ws.Range("G23").Style = "Normal"
$ws.Range("E24").Value = "'1.51%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'7"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.004088"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-4.77%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'7"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001353"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'6.29%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'7"
$ws.Range("G26").Style = "Normal"
$ws.Range("G27").Value = "'7"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'7"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'7"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'7"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'7"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'7"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'7"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'7"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'7"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'7"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'7"
$ws.Range("G37").Style = "Normal"
$ws.Range("D38").Value = "'0.02412"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-2.95%"
$ws.Range("E38").Style = "Normal"
$ws.Range("G38").Value = "'7"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.05147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.08%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'7"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.007747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.21%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'7"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.1302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.76%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'7"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.007603"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'10.49%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'7"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.003304"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'72.64%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'7"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.008493"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.35%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'7"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.3386"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'10.59%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'7"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006588"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.20%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'7"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'7"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.2687"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-38.80%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'7"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.004204"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'35.32%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'7"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'7"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.11%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'7"
$ws.Range("G51").Style = "Normal"
